$wb = $excel.ActiveWorkbook

# --- Rename "Feuil1" -> "Cost Summary_OLD" ---------------------------------
# (Also updates the _xlnm._FilterDatabase defined name that points at it.)
$wsOld = $wb.Worksheets.Item("Feuil1")
$wsOld.Name = "Cost Summary_OLD"

# --- Biomass_Cost: move selection from Z61 to Z27 --------------------------
$wsBiomassCost = $wb.Worksheets.Item("Biomass_Cost")
$wsBiomassCost.Activate()
$wsBiomassCost.Range("Z27").Select()

# --- Cost Analysis-FT: move selection from F29:F37 to F49 ------------------
$wsCostFT = $wb.Worksheets.Item("Cost Analysis-FT")
$wsCostFT.Activate()
$wsCostFT.Range("F49").Select()

# --- Cost Analysis-HEFA: move selection from I25 to F42 --------------------
$wsCostHEFA = $wb.Worksheets.Item("Cost Analysis-HEFA")
$wsCostHEFA.Activate()
$wsCostHEFA.Range("F42").Select()

# --- Kerosene: becomes the active / selected tab ----------------------------
$wsKerosene = $wb.Worksheets.Item("Kerosene")
$wsKerosene.Activate()
$wsKerosene.Range("A1").Select()
